# Atualizações dados 19/07 16h
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 45492.67336805556

# Row -> D,E,F,G values (after). H is set to $newDate for every data row (2..21).
$rows = @{
  2  = @{ D = 1958;  E = 1;  F = 36; G = 17 }
  3  = @{ D = 1963;  E = 2;  F = 33; G = 17 }
  4  = @{ D = 5981;  E = 3;  F = 31; G = 16 }
  5  = @{ D = 1981;  E = 4;  F = 30; G = 17 }
  6  = @{ D = 1955;  E = 5;  F = 30; G = 17 }
  7  = @{ D = 1954;  E = 6;  F = 29; G = 16 }
  8  = @{ D = 2020;  E = 7;  F = 29; G = 16 }
  9  = @{ D = 1967;  E = 8;  F = 25; G = 16 }
  10 = @{ D = 1974;  E = 9;  F = 23; G = 17 }
  11 = @{ D = 1999;  E = 10; F = 22; G = 16 }
  12 = @{ D = 1977;  E = 11; F = 22; G = 16 }
  13 = @{ D = 1980;  E = 12; F = 20; G = 15 }
  14 = @{ D = 1966;  E = 13; F = 19; G = 13 }
  15 = @{ D = 1984;  E = 14; F = 17; G = 15 }
  16 = @{ D = 49202; E = 15; F = 17; G = 16 }
  17 = @{ D = 1962;  E = 16; F = 15; G = 17 }
  18 = @{ D = 1957;  E = 17; F = 15; G = 17 }
  19 = @{ D = 5926;  E = 18; F = 11; G = 15 }
  20 = @{ D = 7314;  E = 19; F = 11; G = 17 }
  21 = @{ D = 1961;  E = 20; F = 8;  G = 16 }
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  $ws.Range("D$r").Value = $vals.D
  $ws.Range("E$r").Value = $vals.E
  $ws.Range("F$r").Value = $vals.F
  $ws.Range("G$r").Value = $vals.G
  $ws.Range("H$r").Value = $newDate
}
